$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.943.61"
Set-TextValue $ws.Range("E2") "  +1.09%  "
Set-TextValue $ws.Range("D3") "1.632.20"
Set-TextValue $ws.Range("E3") "  +2.05%  "
Set-TextValue $ws.Range("E4") "  -0.09%  "
Set-TextValue $ws.Range("D5") "214.74"
Set-TextValue $ws.Range("E5") "  +1.14%  "
Set-TextValue $ws.Range("D6") "0.521"
Set-TextValue $ws.Range("E6") "  +1.24%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.11%  "
Set-TextValue $ws.Range("D8") "29.66"
Set-TextValue $ws.Range("E8") "  +10.40%  "
Set-TextValue $ws.Range("D9") "0.260"
Set-TextValue $ws.Range("E9") "  +3.53%  "
Set-TextValue $ws.Range("E10") "  +2.38%  "
Set-TextValue $ws.Range("E11") "  +0.74%  "
Set-TextValue $ws.Range("D12") "1.863.56"
Set-TextValue $ws.Range("E12") "  +1.98%  "
Set-TextValue $ws.Range("D13") "1.630.87"
Set-TextValue $ws.Range("E13") "  +1.86%  "
Set-TextValue $ws.Range("D14") "0.572"
Set-TextValue $ws.Range("E14") "  +6.49%  "
Set-TextValue $ws.Range("D15") "9.46"
Set-TextValue $ws.Range("E15") "  +24.35%  "
Set-TextValue $ws.Range("D16") "3.89"
Set-TextValue $ws.Range("E16") "  +4.36%  "
Set-TextValue $ws.Range("D17") "29.942.62"
Set-TextValue $ws.Range("E17") "  +1.09%  "
Set-TextValue $ws.Range("D18") "64.86"
Set-TextValue $ws.Range("E18") "  +1.41%  "
Set-TextValue $ws.Range("D19") "249.64"
Set-TextValue $ws.Range("E19") "  +3.37%  "
Set-TextValue $ws.Range("D20") "0.0₃0706"
Set-TextValue $ws.Range("E20") "  +1.87%  "
Set-TextValue $ws.Range("E21") "  -0.04%  "
Set-TextValue $ws.Range("D22") "4.16"
Set-TextValue $ws.Range("E22") "  +4.81%  "
Set-TextValue $ws.Range("D23") "9.62"
Set-TextValue $ws.Range("E23") "  +4.14%  "
Set-TextValue $ws.Range("E24") "  +0.97%  "
Set-TextValue $ws.Range("D25") "159.70"
Set-TextValue $ws.Range("E25") "  +3.39%  "
Set-TextValue $ws.Range("D26") "15.72"
Set-TextValue $ws.Range("E26") "  +2.38%  "
Set-TextValue $ws.Range("E27") "  +2.39%  "
Set-TextValue $ws.Range("D28") "6.61"
Set-TextValue $ws.Range("E28") "  +3.49%  "
Set-TextValue $ws.Range("E29") "  -0.07%  "
Set-TextValue $ws.Range("D30") "0.0490"
Set-TextValue $ws.Range("E30") "  +2.59%  "
Set-TextValue $ws.Range("E31") "  +6.75%  "
Set-TextValue $ws.Range("D32") "3.37"
Set-TextValue $ws.Range("E32") "  +4.79%  "
Set-TextValue $ws.Range("E33") "  +2.01%  "
Set-TextValue $ws.Range("D34") "1.429.92"
Set-TextValue $ws.Range("E34") "  +0.33%  "
Set-TextValue $ws.Range("E35") "  +6.91%  "
Set-TextValue $ws.Range("E36") "  +1.15%  "
Set-TextValue $ws.Range("D37") "2.87"
Set-TextValue $ws.Range("E37") "  -0.29%  "
Set-TextValue $ws.Range("B38") "VeChain"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.0172"
Set-TextValue $ws.Range("E38") "  +3.32%  "
Set-TextValue $ws.Range("B39") "HuobiToken"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D39") "2.29"
Set-TextValue $ws.Range("E39") "  -0.22%  "
Set-TextValue $ws.Range("D40") "0.555"
Set-TextValue $ws.Range("E40") "  +2.35%  "
Set-TextValue $ws.Range("D41") "73.58"
Set-TextValue $ws.Range("E41") "  +11.83%  "
Set-TextValue $ws.Range("E42") "  +1.18%  "
Set-TextValue $ws.Range("D43") "0.832"
Set-TextValue $ws.Range("E43") "  +3.34%  "
Set-TextValue $ws.Range("D44") "55.05"
Set-TextValue $ws.Range("E44") "  +1.04%  "
Set-TextValue $ws.Range("E45") "  +0.30%  "
Set-TextValue $ws.Range("E46") "  +5.21%  "
Set-TextValue $ws.Range("E47") "  -0.05%  "
Set-TextValue $ws.Range("D48") "5.48"
Set-TextValue $ws.Range("E48") "  +3.10%  "
Set-TextValue $ws.Range("D49") "1.770.39"
Set-TextValue $ws.Range("E49") "  +1.91%  "
Set-TextValue $ws.Range("D50") "89.88"
Set-TextValue $ws.Range("E50") "  +4.46%  "
Set-TextValue $ws.Range("E51") "  +1.87%  "
